# Nexial command catalog update ("#system" sheet):
#   - [json] add new function `storeKeys(json,jsonpath,var)` (alphabetically
#     between storeCount and storeValue) -> json range grows from M2:M17 to M2:M18
#   - [target] drop the standalone "text" category from the category list in
#     column A (A2:A31 -> A2:A30); the lone "text" entry (spellCheck) that used
#     to live in its own column Y is merged away and every block that followed
#     it (web, webalert, webcookie, ws, ws.async, xml) shifts one column left
#     (Z->Y, AA->Z, AB->AA, AC->AB, AD->AC, AE->AD).

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("#system")

# ---------------------------------------------------------------------------
# 1) [json] insert "storeKeys(json,jsonpath,var)" at M16, pushing the two
#    existing entries (storeValue, storeValues) down by one row.
# ---------------------------------------------------------------------------
$ws.Range("M17:M18").Value = $ws.Range("M16:M17").Value2
$ws.Range("M16").Value = "storeKeys(json,jsonpath,var)"

# ---------------------------------------------------------------------------
# 2) [target] remove the "text" entry at A25, pulling everything below
#    (web..xml) up by one row; the list now ends at A30 instead of A31.
# ---------------------------------------------------------------------------
$ws.Range("A25:A30").Value = $ws.Range("A26:A31").Value2
$ws.Range("A31").ClearContents()

# ---------------------------------------------------------------------------
# 3) Shift the web/webalert/webcookie/ws/ws.async/xml command blocks one
#    column to the left (Z:AE -> Y:AD), overwriting the now-redundant
#    standalone "text" column Y, and clear the vacated column AE.
# ---------------------------------------------------------------------------
$ws.Range("Y1:AD129").Value = $ws.Range("Z1:AE129").Value2
$ws.Range("AE1:AE129").ClearContents()

# ---------------------------------------------------------------------------
# 4) Update the defined names (named ranges) to match the new layout.
# ---------------------------------------------------------------------------
function Set-SystemName($name, $ref) {
    foreach ($n in $wb.Names) {
        if ($n.Name -eq $name) {
            $n.RefersTo = $ref
        }
    }
}

Set-SystemName "json"      "='#system'!`$M`$2:`$M`$18"
Set-SystemName "target"    "='#system'!`$A`$2:`$A`$30"
Set-SystemName "web"       "='#system'!`$Y`$2:`$Y`$129"
Set-SystemName "webalert"  "='#system'!`$Z`$2:`$Z`$8"
Set-SystemName "webcookie" "='#system'!`$AA`$2:`$AA`$8"
Set-SystemName "ws"        "='#system'!`$AB`$2:`$AB`$17"
Set-SystemName "ws.async"  "='#system'!`$AC`$2:`$AC`$8"
Set-SystemName "xml"       "='#system'!`$AD`$2:`$AD`$27"
